$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently ends with a merged totals row at row 18 (B18:G18 merged,
# "Totaal prijs:" + H18 = SUM(H3:H17)). We need to:
#   1. Insert a new line-item row (new row 18) for "Motorbracket" / "PLA".
#   2. Push the totals row down to row 19, updating its SUM range.
#   3. Keep the rest of the sheet (headers, hyperlinks, etc.) untouched.
# ---------------------------------------------------------------------------

# Unmerge the old totals row so we can move/recreate its cells individually.
$null = $ws.Range("B18:G18").UnMerge()

# Copy the formatting (borders/alignment/number format) of the old totals
# row down to row 19 first, so row 19 ends up with the exact same look the
# totals row always had.
$ws.Range("B18:H18").Copy($ws.Range("B19"))

# Build the new totals row (row 19): label + merged cell, and the updated
# SUM formula that now covers the new line item in row 18.
$null = $ws.Range("B19:G19").Merge()
$ws.Range("B19").Value = "Totaal prijs:"
$ws.Range("H19").Formula = "=SUM(H3:H18)"

# Now turn the old row 18 into the new "Motorbracket" line item. Start by
# clearing whatever leftover content/format it still has, then paste the
# formatting used by the rest of the line-item rows (row 17) onto it so the
# borders/alignment match the rest of the table.
$ws.Range("B18:H18").ClearContents()
$ws.Range("B17:H17").Copy()
$ws.Range("B18").PasteSpecial(-4122)

$ws.Range("B18").Value = 16
$ws.Range("C18").Value = "Motorbracket"
$ws.Range("D18").Value = "PLA"
$ws.Range("E18").Value = "Nieuw"
$ws.Range("F18").Value = 0.34
$ws.Range("G18").Value = 2
$ws.Range("H18").Formula = "=F18*G18"

# Update the selection to match the saved cursor position in the edited file.
$null = $ws.Range("F21").Select()
